# Agregan nuevas entidades de salud en documento y en reporte
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 25-28: copy formatting from the last existing data row (24) ---
$ws.Range("A24:I24").Copy()
$ws.Range("A25:I28").PasteSpecial(-4122)

# Row 25's D cell keeps the "code" style (s=5) instead of the normal column
# style, matching the source workbook's data (copy formatting from B25).
$ws.Range("B25").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# --- Populate the "razonEntidad" (D) column first, then "concepto" (B) ---
# column, then the rest, so new shared-string entries land in the exact
# order the original authoring session produced them.
$ws.Range("D25").Value = "MEDIMAS"
$ws.Range("D26").Value = "COOMEVA"
$ws.Range("D27").Value = "ASOCIACION MUTUAL EMSSANAR"
$ws.Range("D28").Value = "CAFESALUD"

$ws.Range("B25").Value = "COD1"
$ws.Range("B26").Value = "COD2"
$ws.Range("B27").Value = "COD3"
$ws.Range("B28").Value = "COD4"

$ws.Range("A25").Value = 901097473
$ws.Range("C25").Value = "contribuciones"
$ws.Range("E25").Value = "SALUD"
$ws.Range("F25").Value = "A010102002"
$ws.Range("G25").Value = 73
$ws.Range("H25").Value = 23001010102
$ws.Range("I25").Value = "EMP. PRIV. PROMOTORAS DE SALUD"

$ws.Range("A26").Value = 805000427
$ws.Range("C26").Value = "contribuciones"
$ws.Range("E26").Value = "SALUD"
$ws.Range("F26").Value = "A010102002"
$ws.Range("G26").Value = 73
$ws.Range("H26").Value = 23001010102
$ws.Range("I26").Value = "EMP. PRIV. PROMOTORAS DE SALUD"

$ws.Range("A27").Value = 814000337
$ws.Range("C27").Value = "contribuciones"
$ws.Range("E27").Value = "SALUD"
$ws.Range("F27").Value = "A010102002"
$ws.Range("G27").Value = 73
$ws.Range("H27").Value = 23001010102
$ws.Range("I27").Value = "EMP. PRIV. PROMOTORAS DE SALUD"

$ws.Range("A28").Value = 800140949
$ws.Range("C28").Value = "contribuciones"
$ws.Range("E28").Value = "SALUD"
$ws.Range("F28").Value = "A010102002"
$ws.Range("G28").Value = 73
$ws.Range("H28").Value = 23001010102
$ws.Range("I28").Value = "EMP. PRIV. PROMOTORAS DE SALUD"

# --- Workbook-level "_FilterDatabase" hidden defined name over A1:A27 ---
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "='datos Entidades'!`$A`$1:`$A`$27")
$fd.Visible = $false

# --- Final selection left on B28, as in the saved workbook ---
$null = $ws.Range("B28").Select()

Write-Output "done"
